# Addition of some stadistics and a better creation of the variables
#
# Adds two new worksheets ("Partida 5" and "Partida 6") at the end of the
# workbook, each with the same "T","V","A" header as the existing sheets,
# fills them with their round data, and leaves "Partida 6" as the active
# (selected) sheet/tab, matching the target workbook state.

$wb = $excel.ActiveWorkbook

$sheet5Data = @(
  @(10,-5,-5),
  @(10,10,-5),
  @(10,15,-5),
  @(-5,10,20),
  @(-5,10,25),
  @(-5,15,-10),
  @(-10,-5,15),
  @(25,25,-5),
  @(-5,20,20),
  @(-10,20,-5),
  @(20,-5,25),
  @(-5,-5,15),
  @(-5,25,35),
  @(-5,-5,45),
  @(35,-10,35),
  @(30,40,-5),
  @(15,-5,40),
  @(-10,45,20),
  @(-5,30,20),
  @(-10,15,30),
  @(30,-10,25),
  @(-5,25,15),
  @(-5,-5,20),
  @(-5,25,15),
  @(-10,15,-5),
  @(-5,20,10),
  @(10,-5,10)
)

$sheet6Data = @(
  @(15,-5,10),
  @(-5,10,15),
  @(15,10,-5),
  @(10,-5,10),
  @(20,-5,15),
  @(-5,15,-5),
  @(20,-5,25),
  @(20,-10,-5),
  @(20,30,-5),
  @(20,-10,-5),
  @(25,5,25),
  @(-5,-10,40),
  @(25,-5,25),
  @(25,-5,30),
  @(35,-5,20),
  @(-5,35,25),
  @(30,40,-10),
  @(30,-5,-5),
  @(30,-5,30),
  @(20,35,-5),
  @(-5,30,35),
  @(25,-5,-15),
  @(-5,25,15),
  @(-5,20,-5),
  @(25,-5,20),
  @(-5,-10,-5),
  @(20,15,-5),
  @(20,-5,15),
  @(15,-5,10)
)

# --- Partida 5 ---------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "Partida 5"

$ws5.Cells.Item(1,1).Value = "T"
$ws5.Cells.Item(1,2).Value = "V"
$ws5.Cells.Item(1,3).Value = "A"

$r = 2
foreach ($row in $sheet5Data) {
  $ws5.Cells.Item($r,1).Value = $row[0]
  $ws5.Cells.Item($r,2).Value = $row[1]
  $ws5.Cells.Item($r,3).Value = $row[2]
  $r = $r + 1
}

[void]$ws5.Range("M26").Select()

# --- Partida 6 ---------------------------------------------------------
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws5)
$ws6.Name = "Partida 6"

$ws6.Cells.Item(1,1).Value = "T"
$ws6.Cells.Item(1,2).Value = "V"
$ws6.Cells.Item(1,3).Value = "A"

$r = 2
foreach ($row in $sheet6Data) {
  $ws6.Cells.Item($r,1).Value = $row[0]
  $ws6.Cells.Item($r,2).Value = $row[1]
  $ws6.Cells.Item($r,3).Value = $row[2]
  $r = $r + 1
}

[void]$ws6.Range("G31").Select()
